$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The four records currently on rows 45-48 are being re-ordered (a cyclic
# shift): the record that was on row 46 moves up to row 45, the one on row
# 47 moves to row 46, the one on row 48 moves to row 47, and the one that
# was on row 45 moves down to row 48. Row 46's "Aktivitet" value (column M)
# travels together with its record.
#
# Copy/PasteSpecial (rather than reading/writing .Value2 directly) is used
# so that text that merely looks like a date or time (e.g. "2023-08-22",
# "09:36") is preserved verbatim instead of being re-interpreted as a
# date/time serial number by the assignment.

# 1) Stash the original row 45 in a scratch row far outside the used range.
$ws.Range("A45:AY45").Copy() | Out-Null
$ws.Range("A1000:AY1000").PasteSpecial() | Out-Null

# 2) row46 -> row45 (row46 has a value in column M; row45 currently has none)
$ws.Range("M45").Clear() | Out-Null
$ws.Range("A46:AY46").Copy() | Out-Null
$ws.Range("A45:AY45").PasteSpecial() | Out-Null

# 3) row47 -> row46 (row47 has no M value, so clear row46's M first)
$ws.Range("M46").Clear() | Out-Null
$ws.Range("A47:AY47").Copy() | Out-Null
$ws.Range("A46:AY46").PasteSpecial() | Out-Null

# 4) row48 -> row47 (row48 has no M value, so clear row47's M first)
$ws.Range("M47").Clear() | Out-Null
$ws.Range("A48:AY48").Copy() | Out-Null
$ws.Range("A47:AY47").PasteSpecial() | Out-Null

# 5) scratch (original row45) -> row48 (original row45 has no M value)
$ws.Range("M48").Clear() | Out-Null
$ws.Range("A1000:AY1000").Copy() | Out-Null
$ws.Range("A48:AY48").PasteSpecial() | Out-Null

# 6) Clean up the scratch row.
$ws.Range("A1000:AY1000").Clear() | Out-Null
